$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.16, 3911.892387694, 37.487058824, 3874.405328872, 25.440260226, 1585.105019658, 23.140807668, 1561.96421199, 26.67706112333333),
    @(0.17, 4041.011614314, 36.964767414, 4004.0468469, 25.44409872933333, 1608.56998527, 22.774044308, 1585.795940964, 26.68560294533333),
    @(0.18, 4312.789545273999, 36.30853267600001, 4276.481012597999, 25.30718837866667, 1696.902866872, 22.475940104, 1674.426926766, 26.55429984466667),
    @(0.19, 4427.371668214, 35.672811882, 4391.698856334, 25.31023158866667, 1708.110787172, 22.116897024, 1685.993890148, 26.56685283933334),
    @(0.2, 4719.73852585, 35.270357802, 4684.468168048, 25.32196972466667, 1798.145568474, 21.926013394, 1776.21955508, 26.58711745933334)
)

$startRow = 17
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
